# Finalized definition of virtual environment
# Update calibration parameter values (and their derived autocorrelation
# array strings) across the AR / SETAR / GARCH / TARCH / AR_TARCH sheets.

$wb = $excel.ActiveWorkbook

# --- AR sheet ---
$ws = $wb.Worksheets.Item("AR")
$ws.Range("B2").Value = -0.01403177878409627
$ws.Range("B3").Value = 0.8322473199251882
$ws.Range("B4").Value = 34.65089605793457
$ws.Range("B5").Value = "[1.0, 0.30237599728096964, 0.24218465485834614, 0.27410815183423765, 0.25513874014837684, 0.4669657394855553, 0.24733522757548726, 0.1944752928840054, 0.23746590281237834, 0.19977431858980044, 0.1945605552084176, 0.13098935671650389, 0.1383911102598105, 0.19585270706277116, 0.1955989684696984, 0.1531427091844605, 0.15607934026067333, 0.1394063095344178, 0.14903665562575663, 0.16778586445993193]"

# --- SETAR sheet ---
$ws = $wb.Worksheets.Item("SETAR")
$ws.Range("B3").Value = 0.9434219495569189
$ws.Range("B4").Value = -1.576957632985456
$ws.Range("B5").Value = 0.7790702448451842
$ws.Range("B6").Value = 27.26141873251629
$ws.Range("B7").Value = 1.438973502798857
$ws.Range("B8").Value = 0.780899905860917
$ws.Range("B9").Value = 24.3998919420025
$ws.Range("B10").Value = "[1.0, 0.2892229777575205, 0.2425376463818646, 0.2758507447187141, 0.2617769154591766, 0.34882327875039787, 0.2128903348391278, 0.17993374273571347, 0.1968467281166908, 0.18774935300069656, 0.19100315638568155, 0.14099496173922071, 0.14030667199996075, 0.1576689547271237, 0.14672875698233195, 0.1301287465339835, 0.15562221153768904, 0.09649613656004367, 0.125133412563086, 0.15863811952239476]"

# --- GARCH sheet ---
$ws = $wb.Worksheets.Item("GARCH")
$ws.Range("B2").Value = 0.006651707337102347
$ws.Range("B3").Value = 0.2998353563667131
$ws.Range("B4").Value = 0.1380453368187803
$ws.Range("B5").Value = 0.8619546631832927

# --- TARCH sheet ---
$ws = $wb.Worksheets.Item("TARCH")
$ws.Range("B2").Value = -0.005992606524840025
$ws.Range("B3").Value = 0.300788763197276
$ws.Range("B4").Value = 0.1311433819625628
$ws.Range("B5").Value = 0.8618197969154329
$ws.Range("B6").Value = "[1.0, 0.01375426401571509, -0.034814428767388515, -0.00870144925032259, -0.03692035743757103, 0.29950961618291466, -0.009779334061821654, -0.01148855682729653, 0.019356383993622273, -0.033004067880131265, -0.058934078787601166, -0.07741875383981178, -0.06043646440267581, -0.006740112244389282, 0.014401611512009416, -0.02058280047067888, -0.024740679142091535, -0.051228246314968036, -0.0628300690287674, -0.04314532787565816]"
$ws.Range("B7").Value = 0.01407365119640938

# --- AR_TARCH sheet ---
$ws = $wb.Worksheets.Item("AR_TARCH")
$ws.Range("B2").Value = -0.02505925245946669
$ws.Range("B3").Value = 0.293231811634583
$ws.Range("B4").Value = 0.1421410966524452
$ws.Range("B5").Value = 0.8599586795188392
$ws.Range("B6").Value = "[1.0, 0.017383593523458682, -0.022132188568854946, 0.013217038015203048, -0.026054202167705275, 0.23495178609310788, -0.0081535301041223, -0.010042423485942767, 0.004483229686143032, -0.044259425602285865, -0.029461128600988876, -0.0690679571503173, -0.05341582454190929, -0.008774250307792702, 0.013899969311116539, -0.03930143259795254, -0.019268778574064872, -0.03586988819557108, -0.04565582640326065, -0.02564643797459133]"
$ws.Range("B7").Value = -0.004199552308136415
$ws.Range("B9").Value = 0.8031892808210191
